# Add three more demo users (user2, user3, user4) below the existing
# "user1" row, each paired with the same john/demo credentials, so the
# login-data sheet supports the new parallel scenario runners.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 keeps its username/password pair but the login id becomes "user1".
$ws.Cells.Item(2, 1).Value = "user1"
$ws.Cells.Item(2, 2).Value = "john"
$ws.Cells.Item(2, 3).Value = "demo"

# New rows 3-5: user2..user4, same username/password as row 2.
for ($i = 2; $i -le 4; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = "user$i"
    $ws.Cells.Item($row, 2).Value = "john"
    $ws.Cells.Item($row, 3).Value = "demo"
}

$ws.Range("C2:C5").Select() | Out-Null
